$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$siteA = "四方坪站"
$siteB = "高岭站"

# Row 44 (F44 pre-existed as an empty placeholder cell with the wrong
# number format, so restore the column's integer format explicitly)
$ws.Cells.Item(44, 1).Value = 46075
$ws.Cells.Item(44, 2).Value = $siteA
$ws.Cells.Item(44, 3).Value = 7053.32
$ws.Cells.Item(44, 4).Value = 6673.99
$ws.Cells.Item(44, 5).Value = 2629.38
$ws.Cells.Item(44, 6).NumberFormat = '0_);[Red]\(0\)'
$ws.Cells.Item(44, 6).Value = 277

# Row 45
$ws.Cells.Item(45, 1).Value = 46075
$ws.Cells.Item(45, 2).Value = $siteB
$ws.Cells.Item(45, 3).Value = 1446.01
$ws.Cells.Item(45, 4).Value = 1315.6
$ws.Cells.Item(45, 5).Value = 407.29
$ws.Cells.Item(45, 6).Value = 49

# Row 46
$ws.Cells.Item(46, 1).Value = 46076
$ws.Cells.Item(46, 2).Value = $siteA
$ws.Cells.Item(46, 3).Value = 7385.61
$ws.Cells.Item(46, 4).Value = 7003.54
$ws.Cells.Item(46, 5).Value = 2733.92
$ws.Cells.Item(46, 6).Value = 320

# Row 47
$ws.Cells.Item(47, 1).Value = 46076
$ws.Cells.Item(47, 2).Value = $siteB
$ws.Cells.Item(47, 3).Value = 2221.13
$ws.Cells.Item(47, 4).Value = 2041.25
$ws.Cells.Item(47, 5).Value = 617.02
$ws.Cells.Item(47, 6).Value = 70

# Row 48
$ws.Cells.Item(48, 1).Value = 46077
$ws.Cells.Item(48, 2).Value = $siteA
$ws.Cells.Item(48, 3).Value = 9399.94
$ws.Cells.Item(48, 4).Value = 8724.5300000000007
$ws.Cells.Item(48, 5).Value = 3458.13
$ws.Cells.Item(48, 6).Value = 383

# Row 49
$ws.Cells.Item(49, 1).Value = 46077
$ws.Cells.Item(49, 2).Value = $siteB
$ws.Cells.Item(49, 3).Value = 2913.86
$ws.Cells.Item(49, 4).Value = 2848.88
$ws.Cells.Item(49, 5).Value = 791.49
$ws.Cells.Item(49, 6).Value = 108

# Update view: scroll position and selection
$ws.Range("H50").Select()
$excel.ActiveWindow.ScrollRow = 40
